$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value looks like a plain number; force Text format
# so Excel stores the literal string (matching the source data which keeps
# prices as text), not an auto-converted number.
$textCells = @("D5", "D8", "D9", "D11", "D16", "D17", "D18", "D23", "D25", "D30", "D33", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped from the latest cryptos snapshot.
$ws.Range("D2").Value = "28.299.99"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "1.550.36"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("D5").Value = "209.39"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("E6").Value = "  -1.37%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "23.64"
$ws.Range("E8").Value = "  -1.83%  "
$ws.Range("D9").Value = "0.241"
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").Value = "0.0890"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "1.772.41"
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").Value = "1.546.85"
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").Value = "28.288.99"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "0.508"
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("D17").Value = "60.54"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").Value = "226.84"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").Value = "8.82"
$ws.Range("E23").Value = "  -3.18%  "
$ws.Range("E24").Value = "  -6.06%  "
$ws.Range("D25").Value = "149.29"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("E26").Value = "  -1.73%  "
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D30").Value = "0.0465"
$ws.Range("E30").Value = "  -3.94%  "
$ws.Range("E31").Value = "  -4.65%  "
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("D33").Value = "3.04"
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("D34").Value = "1.381.07"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("D38").Value = "2.59"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "1.90"
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("B42").Value = "ImmutableX"
$ws.Range("C42").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D42").Value = "0.507"
$ws.Range("E42").Value = "  -3.24%  "
$ws.Range("D43").Value = "0.774"
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("D44").Value = "0.0464"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("D45").Value = "5.37"
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("D46").Value = "61.76"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.685.66"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "0.906"
$ws.Range("E48").Value = "  -6.46%  "
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").Value = "42.05"
$ws.Range("E50").Value = "  +5.89%  "
$ws.Range("E51").Value = "  +0.37%  "
